$wb = $excel.ActiveWorkbook

# --- Refresh the panel query timestamps on the "data" sheet ------------
$data = $wb.Worksheets.Item("data")

$newDataTimestamps = @(
"2021-10-05 14:20:37.335398","2021-10-05 14:20:37.335407","2021-10-05 14:20:37.335410","2021-10-05 14:20:37.335412","2021-10-05 14:20:37.335415","2021-10-05 14:20:37.335418","2021-10-05 14:20:37.335420","2021-10-05 14:20:37.335423","2021-10-05 14:20:37.335426","2021-10-05 14:20:37.335428","2021-10-05 14:20:37.335431","2021-10-05 14:20:37.335433","2021-10-05 14:20:37.335436","2021-10-05 14:20:37.335438","2021-10-05 14:20:37.335441","2021-10-05 14:20:37.335443","2021-10-05 14:20:37.335446","2021-10-05 14:20:37.335449","2021-10-05 14:20:37.335451","2021-10-05 14:20:37.335454","2021-10-05 14:20:37.335456","2021-10-05 14:20:37.335459","2021-10-05 14:20:37.335461","2021-10-05 14:20:37.335464","2021-10-05 14:20:37.335467","2021-10-05 14:20:37.335469","2021-10-05 14:20:37.335472","2021-10-05 14:20:37.335474","2021-10-05 14:20:37.335477","2021-10-05 14:20:37.335479","2021-10-05 14:20:37.335482","2021-10-05 14:20:37.335484","2021-10-05 14:20:37.335487","2021-10-05 14:20:37.335490","2021-10-05 14:20:37.335492","2021-10-05 14:20:37.335495","2021-10-05 14:20:37.335497","2021-10-05 14:20:37.335500","2021-10-05 14:20:37.335502","2021-10-05 14:20:37.335505","2021-10-05 14:20:37.335508","2021-10-05 14:20:37.335511","2021-10-05 14:20:37.335513","2021-10-05 14:20:37.335516","2021-10-05 14:20:37.335518","2021-10-05 14:20:37.335521","2021-10-05 14:20:37.335523","2021-10-05 14:20:37.335526","2021-10-05 14:20:37.335528","2021-10-05 14:20:37.335530","2021-10-05 14:20:37.335533","2021-10-05 14:20:37.335535","2021-10-05 14:20:37.335538","2021-10-05 14:20:37.335541","2021-10-05 14:20:37.335543","2021-10-05 14:20:37.335546","2021-10-05 14:20:37.335549","2021-10-05 14:20:37.335551","2021-10-05 14:20:37.335554","2021-10-05 14:20:37.335557","2021-10-05 14:20:37.335559","2021-10-05 14:20:37.335562","2021-10-05 14:20:37.335564","2021-10-05 14:20:37.335567","2021-10-05 14:20:37.335570","2021-10-05 14:20:37.335573","2021-10-05 14:20:37.335576","2021-10-05 14:20:37.335578","2021-10-05 14:20:37.335581","2021-10-05 14:20:37.335583","2021-10-05 14:20:37.335586","2021-10-05 14:20:37.335588","2021-10-05 14:20:37.335591","2021-10-05 14:20:37.335594","2021-10-05 14:20:37.335596","2021-10-05 14:20:37.335599","2021-10-05 14:20:37.335603","2021-10-05 14:20:37.335607","2021-10-05 14:20:37.335609","2021-10-05 14:20:37.335612","2021-10-05 14:20:37.335614","2021-10-05 14:20:37.335617","2021-10-05 14:20:37.335619","2021-10-05 14:20:37.335622","2021-10-05 14:20:37.335625","2021-10-05 14:20:37.335628","2021-10-05 14:20:37.335630","2021-10-05 14:20:37.335633","2021-10-05 14:20:37.335635"
)

for ($i = 0; $i -lt $newDataTimestamps.Length; $i++) {
    $row = $i + 2
    $data.Cells.Item($row, 6).Value = $newDataTimestamps[$i]
}

# --- Add the new "metadata" sheet, right after "data" ------------------
$meta = $wb.Worksheets.Add($null, $data)
$meta.Name = "metadata"

$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

$headerRng = $meta.Range("B1:G1")
$headerRng.Font.Bold = $true
$headerRng.HorizontalAlignment = -4108
$headerRng.VerticalAlignment = -4160
$headerRng.Borders.LineStyle = 1
$headerRng.Borders.Weight = 2

$a2 = $meta.Range("A2")
$a2.Value = 0
$a2.Font.Bold = $true
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160
$a2.Borders.LineStyle = 1
$a2.Borders.Weight = 2

$meta.Range("B2").Value = "Haematological malignancies for rare disease"
$meta.Range("C2").Value = 407
$meta.Range("D2").Value = "1.5"
$meta.Range("E2").Value = "2021-10-01T08:48:00.201378Z"
$meta.Range("F2").Value = "2021-10-05 14:20:37.331983"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/407/?format=json"
